$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "[-, -, -, 'MEC-1NA-Tec. Mat. Não Metal.']"

$ws.Range("C19").Value = "-"
$ws.Range("E19").Value = "[-, -, -, 'MEC-1NA-Tec. Mat. Não Metal.']"

$ws.Range("E20").Value = "[-, -, -, 'MEC-1NA-Tec. Mat. Não Metal.']"

$ws.Range("E21").Value = "[-, -, -, 'MEC-1NA-Tec. Mat. Não Metal.']"
